$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 43: Tuesday, July 2 2024 (serial 45475) ---
$ws.Range("A42:C42").Copy() | Out-Null
$ws.Range("A43:C43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E42").Copy() | Out-Null
$ws.Range("E43").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A43").Value = 45475
$ws.Range("B43").Value = "T"
$ws.Range("C43").Value = 4
$ws.Range("E43").Value = "finalizing studentgrades_prof, cleaning out irrelevant prof columns"
$ws.Rows.Item(43).RowHeight = 28.5

# --- Row 44: Wednesday, July 3 2024 (serial 45476) ---
$ws.Range("A42:C42").Copy() | Out-Null
$ws.Range("A44:C44").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E42").Copy() | Out-Null
$ws.Range("E44").PasteSpecial(-4122) | Out-Null       # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A44").Value = 45476
$ws.Range("B44").Value = "W"
$ws.Range("C44").Value = 4
$ws.Range("E44").Value = "running GBM with different configurations, fixing prof column removal, adding recursive removal to columns to improve prediction"
$ws.Rows.Item(44).RowHeight = 42.75

# --- View: scroll window down, select the new last cell ---
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("E44").Select() | Out-Null
